$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core data fix ("Fixed monitoredData endpoint glitch") ---
# B2 (monitored_goods) held a stale/incorrect timestamp; correct it.
$ws.Range("B2").Value = 45078

# --- Cosmetic state captured by the same save (best effort under this host) ---
# Column B was resized to fit its (long, timestamp-formatted) contents.
$ws.Columns("B").ColumnWidth = 17.45

# Cursor/selection was left on G9 when the author saved the workbook.
[void]$ws.Range("G9").Select()

# A confidentiality footer (typical of an emailed/automated report) was added.
$footerText = "`r&1#&`"Calibri`"&7&K000000 Este conteúdo e quaisquer informações anexadas a ele são confidenciais e destinados exclusivamente para uso do indivíduo ou pela entidade a quem estão endereçados. Se você recebeu este email por engano, notifique o administrador"
$ws.PageSetup.CenterFooter = $footerText
